# Apply the updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.789.76"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.886.23"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'238.68"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4759"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "'0.2872"
$ws.Range("E8").Value = "  +5.10%  "
$ws.Range("D9").Value = "'0.06568"
$ws.Range("E9").Value = "  +4.33%  "
$ws.Range("D10").Value = "'18.91"
$ws.Range("E10").Value = "  +15.88%  "
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "'97.43"
$ws.Range("E11").Value = "  +16.07%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.872.88"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'5.124"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "'0.6563"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").Value = "'309.52"
$ws.Range("E16").Value = "  +35.18%  "
$ws.Range("D17").Value = "30.772.33"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "'13.19"
$ws.Range("E18").Value = "  +6.59%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "'0.000007577"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").Value = "2.123.08"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "'0.9986"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'5.128"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("D24").Value = "'6.191"
$ws.Range("E24").Value = "  +5.43%  "
$ws.Range("D25").Value = "'9.294"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'166.58"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'20.23"
$ws.Range("E27").Value = "  +13.44%  "
$ws.Range("D28").Value = "'1.949"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "'0.1077"
$ws.Range("E29").Value = "  +5.60%  "
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'4.160"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").Value = "'3.978"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("D33").Value = "'0.05046"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").Value = "'1.177"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("D35").Value = "'0.7339"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("D36").Value = "'2.711"
$ws.Range("D37").Value = "'0.01947"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").Value = "'2.701"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "'2.078"
$ws.Range("D40").Value = "'0.9047"
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("D41").Value = "'107.88"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").Value = "'0.9987"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.4223"
$ws.Range("E43").Value = "  +4.82%  "
$ws.Range("D44").Value = "'5.642"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "'65.94"
$ws.Range("E45").Value = "  +7.29%  "
$ws.Range("D46").Value = "'7.385"
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1227"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.988"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("D49").Value = "'34.82"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").Value = "'0.05609"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "'0.3859"
$ws.Range("E51").Value = "  +5.69%  "
